# Weekly update: insert the latest week's two price observations
# (Primera / Segunda calidad) for "Apio" at Feria Lagunitas de Puerto Montt,
# right after the existing row 140, pushing the rest of the historical
# rows (old 141-176) down by two rows to 143-178.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows (shifts rows 141:176 down to 143:178,
# copying the formatting of row 141 - e.g. the date style on column D).
$ws.Rows("141:142").Insert()

# --- New row 141: Apio, Primera, week of 2021-12-21 ---
$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44551
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = 100112017
$ws.Cells.Item(141, 7).Value = "Apio"
$ws.Cells.Item(141, 8).Value = "Americana (o)"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 20
$ws.Cells.Item(141, 11).Value = 12000
$ws.Cells.Item(141, 12).Value = 12000
$ws.Cells.Item(141, 13).Value = 12000
$ws.Cells.Item(141, 14).Value = "$/docena de matas"
$ws.Cells.Item(141, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(141, 16).Value = 2000
$ws.Cells.Item(141, 17).Value = 6
$ws.Cells.Item(141, 18).Value = "Hortaliza"

# --- New row 142: Apio, Segunda, week of 2021-12-21 ---
$ws.Cells.Item(142, 1).Value = 4
$ws.Cells.Item(142, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(142, 3).Value = "Los Lagos"
$ws.Cells.Item(142, 4).Value = 44551
$ws.Cells.Item(142, 5).Value = 10
$ws.Cells.Item(142, 6).Value = 100112017
$ws.Cells.Item(142, 7).Value = "Apio"
$ws.Cells.Item(142, 8).Value = "Americana (o)"
$ws.Cells.Item(142, 9).Value = "Segunda"
$ws.Cells.Item(142, 10).Value = 20
$ws.Cells.Item(142, 11).Value = 10000
$ws.Cells.Item(142, 12).Value = 10000
$ws.Cells.Item(142, 13).Value = 10000
$ws.Cells.Item(142, 14).Value = "$/docena de matas"
$ws.Cells.Item(142, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(142, 16).Value = 1667
$ws.Cells.Item(142, 17).Value = 6
$ws.Cells.Item(142, 18).Value = "Hortaliza"
